{"js": "// Replace the placeholder ellipsis \"\u2026\u2026\u2026\u2026.\" (NRP value for the gelar\n// perkara signatory) with the \"${nrp_gp}\" merge-field token.\nconst body = context.document.body;\nconst results = body.search(\"\u2026\u2026\u2026\u2026.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"${nrp_gp}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the placeholder ellipsis used for the NRP value of the gelar\n# perkara signatory with the \"${nrp_gp}\" merge-field token.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u2026\u2026\u2026\u2026.\"\n$find.Replacement.Text = \"`${nrp_gp}\"\n$find.Forward = $true\n$find.Wrap = 1               # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
